$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy format from E1 (header style), then set text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
"2021-10-05 10:51:26.217463",
"2021-10-05 10:51:26.217474",
"2021-10-05 10:51:26.217479",
"2021-10-05 10:51:26.217481",
"2021-10-05 10:51:26.217484",
"2021-10-05 10:51:26.217487",
"2021-10-05 10:51:26.217490",
"2021-10-05 10:51:26.217492",
"2021-10-05 10:51:26.217495",
"2021-10-05 10:51:26.217497",
"2021-10-05 10:51:26.217500",
"2021-10-05 10:51:26.217502",
"2021-10-05 10:51:26.217505",
"2021-10-05 10:51:26.217507",
"2021-10-05 10:51:26.217510",
"2021-10-05 10:51:26.217512",
"2021-10-05 10:51:26.217515",
"2021-10-05 10:51:26.217518",
"2021-10-05 10:51:26.217520",
"2021-10-05 10:51:26.217523",
"2021-10-05 10:51:26.217525",
"2021-10-05 10:51:26.217528",
"2021-10-05 10:51:26.217530",
"2021-10-05 10:51:26.217533",
"2021-10-05 10:51:26.217536",
"2021-10-05 10:51:26.217538",
"2021-10-05 10:51:26.217541",
"2021-10-05 10:51:26.217543",
"2021-10-05 10:51:26.217546",
"2021-10-05 10:51:26.217549",
"2021-10-05 10:51:26.217551",
"2021-10-05 10:51:26.217554",
"2021-10-05 10:51:26.217557",
"2021-10-05 10:51:26.217560",
"2021-10-05 10:51:26.217562",
"2021-10-05 10:51:26.217565",
"2021-10-05 10:51:26.217567",
"2021-10-05 10:51:26.217570",
"2021-10-05 10:51:26.217572",
"2021-10-05 10:51:26.217575",
"2021-10-05 10:51:26.217578",
"2021-10-05 10:51:26.217580",
"2021-10-05 10:51:26.217583",
"2021-10-05 10:51:26.217585",
"2021-10-05 10:51:26.217588",
"2021-10-05 10:51:26.217590",
"2021-10-05 10:51:26.217593",
"2021-10-05 10:51:26.217595",
"2021-10-05 10:51:26.217598",
"2021-10-05 10:51:26.217601",
"2021-10-05 10:51:26.217603",
"2021-10-05 10:51:26.217605",
"2021-10-05 10:51:26.217608",
"2021-10-05 10:51:26.217611",
"2021-10-05 10:51:26.217614",
"2021-10-05 10:51:26.217616",
"2021-10-05 10:51:26.217619",
"2021-10-05 10:51:26.217621",
"2021-10-05 10:51:26.217624",
"2021-10-05 10:51:26.217627",
"2021-10-05 10:51:26.217629",
"2021-10-05 10:51:26.217632",
"2021-10-05 10:51:26.217634",
"2021-10-05 10:51:26.217637",
"2021-10-05 10:51:26.217641",
"2021-10-05 10:51:26.217644",
"2021-10-05 10:51:26.217646",
"2021-10-05 10:51:26.217649",
"2021-10-05 10:51:26.217651",
"2021-10-05 10:51:26.217654",
"2021-10-05 10:51:26.217656",
"2021-10-05 10:51:26.217659",
"2021-10-05 10:51:26.217661",
"2021-10-05 10:51:26.217664",
"2021-10-05 10:51:26.217666",
"2021-10-05 10:51:26.217669",
"2021-10-05 10:51:26.217674",
"2021-10-05 10:51:26.217677",
"2021-10-05 10:51:26.217679",
"2021-10-05 10:51:26.217682",
"2021-10-05 10:51:26.217684",
"2021-10-05 10:51:26.217687",
"2021-10-05 10:51:26.217689",
"2021-10-05 10:51:26.217692",
"2021-10-05 10:51:26.217695",
"2021-10-05 10:51:26.217697",
"2021-10-05 10:51:26.217700"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
